# Applies the v5.0 "data rate" refresh to Sheet1: replaces the busy-slot-fee
# percentages (rows 2-10) with the new rate figures, and appends the new row 11.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# row -> (A, B, C, D) values, written as literal text (matches the source data,
# which stores these numbers as shared strings rather than numeric cells).
$rows = @(
    @("1422.2747796302056", "1441.6346421409246", "1347.4374935524736", "1196.773507777089"),
    @("2550.601144660998", "2528.4729375523307", "1935.8579236571452", "1430.7319839939926"),
    @("2809.413954511729", "3061.0654578163508", "2090.0649073666095", "1449.0582954114527"),
    @("2813.113169943492", "3294.428960483938", "2135.2383084950843", "1457.5876402571057"),
    @("4757.190419899803", "3363.1984119754225", "2149.9146372633663", "1456.1873878143724"),
    @("4800.578296169801", "3398.197313135931", "2152.4542693317403", "1456.6171879785686"),
    @("4814.685623542333", "3401.083382959531", "2150.7968005116554", "1454.8448981742683"),
    @("4817.958036112066", "3415.19314854679", "2152.0108386119728", "1451.27359382712"),
    @("4816.0332673319535", "3424.4133067873195", "2153.7298768129945", "1461.223081997538"),
    @("4816.7730140177555", "3416.4361573564097", "2148.839550928931", "1454.6518318798164")
)

$rowIndex = 2
foreach ($rowValues in $rows) {
    $colIndex = 1
    foreach ($val in $rowValues) {
        $cell = $ws.Cells.Item($rowIndex, $colIndex)
        # Quote the literal so Excel stores it as text ("=""<val>""") instead of
        # auto-converting the numeric-looking string to a Number cell.
        $cell.Formula = '="' + $val + '"'
        $colIndex++
    }
    $rowIndex++
}

# Convert the formulas above into plain text values in-place (Copy + PasteSpecial
# xlPasteValues = -4163) so the saved file has literal shared-string cells -- no
# formulas, no "Text" number-format style -- identical in shape to the original.
$target = $ws.Range("A2:D11")
$target.Copy()
$target.PasteSpecial(-4163)
